$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the E74 checkpoint marker cell (borrado de checkpoints)
$ws.Range("E74").Clear() | Out-Null

# Copy the ID-column (A) formatting down for the new rows so new entries
# match the existing bordered/bold style used by the word-list id column.
# (Source range size must match the destination size exactly: 59 rows.)
$ws.Range("A16:A74").Copy() | Out-Null
$ws.Range("A75:A133").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New vocabulary rows (word list checkpoint update)
$rows = @(
    @{ Row=75; A=171; B='earn'; C='ganar' },
    @{ Row=76; A=172; B='edge'; C='borde' },
    @{ Row=77; A=173; B='empty'; C='vacio'; D='vaciar' },
    @{ Row=78; A=174; B='enclose'; C='incluir' },
    @{ Row=79; A=175; B='encourage'; C='estimular' },
    @{ Row=80; A=176; B='engine'; C='motor' },
    @{ Row=81; A=177; B='fault'; C='falta' },
    @{ Row=82; A=178; B='faulty'; C='con errores' },
    @{ Row=83; A=179; B='fold'; C='doblar' },
    @{ Row=84; A=180; B='frame'; C='marco' },
    @{ Row=85; A=181; B='freeze'; C='congelar' },
    @{ Row=86; A=182; B='guard'; C='guardia'; D='proteger' },
    @{ Row=87; A=183; B='guess'; C='suponer' },
    @{ Row=88; A=184; B='guide'; C='guia'; D='guiar' },
    @{ Row=89; A=185; B='hammer'; C='martillo'; D='martillar' },
    @{ Row=90; A=186; B='handle'; C='manilla' },
    @{ Row=91; A=187; B='harm'; C='dano'; D='danar' },
    @{ Row=92; A=188; B='harmful'; C='danino' },
    @{ Row=93; A=189; B='harmless'; C='inofensivo' },
    @{ Row=94; A=190; B='height'; C='altura' },
    @{ Row=95; A=191; B='hesitate'; C='dudar' },
    @{ Row=96; A=192; B='hide'; C='esconder' },
    @{ Row=97; A=193; B='hinder'; C='impedir' },
    @{ Row=98; A=194; B='hit'; C='golpear' },
    @{ Row=99; A=195; B='hollow'; C='hueco'; D='vacio' },
    @{ Row=100; A=196; B='hook'; C='gancho'; D='enganchar' },
    @{ Row=101; A=197; B='hurry'; C='apuro'; D='apurarse' },
    @{ Row=102; A=198; B='ice'; C='hielo' },
    @{ Row=103; A=199; B='improve'; C='mejorar' },
    @{ Row=104; A=200; B='improvement'; C='mejora' },
    @{ Row=105; A=201; B='Inquire / enquire'; C='averiguar' },
    @{ Row=106; A=202; B='inquiry'; C='investigacion' },
    @{ Row=107; A=203; B='insurance'; C='seguro' },
    @{ Row=108; A=204; B='intend'; C='intentar' },
    @{ Row=109; A=205; B='inward'; C='hacia adentro' },
    @{ Row=110; A=206; B='jump'; C='saltar' },
    @{ Row=111; A=207; B='key'; C='llave'; D='clave' },
    @{ Row=112; A=208; B='knot'; C='nudo' },
    @{ Row=113; A=209; B='lamp'; C='lampara' },
    @{ Row=114; A=210; B='lean'; C='sin grasa' },
    @{ Row=115; A=211; B='leg'; C='pierna' },
    @{ Row=116; A=212; B='lend'; C='prestar' },
    @{ Row=117; A=213; B='lid'; C='tapa' },
    @{ Row=118; A=214; B='load'; C='carga'; D='cargar' },
    @{ Row=119; A=215; B='lock'; C='cerradura'; D='trancar' },
    @{ Row=120; A=216; B='loose'; C='flojo'; D='holgado' },
    @{ Row=121; A=217; B='loosen'; C='aflojar' },
    @{ Row=122; A=218; B='loud'; C='a alto volumen' },
    @{ Row=123; A=219; B='lump'; C='grumo' },
    @{ Row=124; A=220; B='male'; C='masculino'; D='macho' },
    @{ Row=125; A=221; B='manage'; C='gerenciar' },
    @{ Row=126; A=222; B='management'; C='gerencia' },
    @{ Row=127; A=223; B='match'; C='fosforo'; D='combinar' },
    @{ Row=128; A=224; B='melt'; C='derretir' },
    @{ Row=129; A=225; B='mend'; C='remendar' },
    @{ Row=130; A=226; B='mild'; C='suave'; D='leve' },
    @{ Row=131; A=227; B='mistake'; C='error' },
    @{ Row=132; A=228; B='mix'; C='mezclar' },
    @{ Row=133; A=229; B='mixture'; C='mezcla' }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    if ($r.ContainsKey("D")) {
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }
    $ws.Cells.Item($r.Row, 5).Value = 0
}

# Restore view state (scrolled down, last edited cell selected)
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("B116").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 116
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D132").Select() | Out-Null

